# Lesson-19 wordlist: re-sequence vocabulary rows into their new section order
# (honorific-speech words move up to rows 2-10/21-44; the family/season/
#  kanji-topic block 54-120 is reshuffled into new section groupings).
# Row 1 (header) and rows 11-20 / 45-53 keep their original content and position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot every data row's English/Japanese pair up front. The reorder is made
# up of long row-cycles (e.g. 2 -> 21 -> 30 -> ... -> 36 -> 2), so source cells
# must all be read before any destination cell is overwritten.
$origA = @{}
$origB = @{}
for ($r = 2; $r -le 120; $r++) {
    $origA[$r] = $ws.Cells.Item($r, 1).Value2
    $origB[$r] = $ws.Cells.Item($r, 2).Value2
}

# destRow -> srcRow: destRow will receive the pair that used to sit at srcRow
$moves = @(
    ,@(2, 21)
    ,@(3, 22)
    ,@(4, 23)
    ,@(5, 24)
    ,@(6, 25)
    ,@(7, 26)
    ,@(8, 27)
    ,@(9, 28)
    ,@(10, 29)
    ,@(21, 30)
    ,@(22, 31)
    ,@(23, 32)
    ,@(24, 33)
    ,@(25, 34)
    ,@(26, 35)
    ,@(27, 36)
    ,@(28, 37)
    ,@(29, 38)
    ,@(30, 39)
    ,@(31, 40)
    ,@(32, 41)
    ,@(33, 42)
    ,@(34, 43)
    ,@(35, 44)
    ,@(36, 2)
    ,@(37, 3)
    ,@(38, 4)
    ,@(39, 5)
    ,@(40, 6)
    ,@(41, 7)
    ,@(42, 8)
    ,@(43, 9)
    ,@(44, 10)
    ,@(54, 78)
    ,@(55, 79)
    ,@(56, 80)
    ,@(57, 81)
    ,@(58, 82)
    ,@(59, 83)
    ,@(60, 84)
    ,@(61, 85)
    ,@(62, 86)
    ,@(63, 87)
    ,@(64, 88)
    ,@(65, 108)
    ,@(66, 109)
    ,@(67, 110)
    ,@(68, 111)
    ,@(69, 112)
    ,@(70, 113)
    ,@(71, 114)
    ,@(72, 115)
    ,@(73, 116)
    ,@(74, 117)
    ,@(75, 118)
    ,@(76, 119)
    ,@(77, 120)
    ,@(78, 54)
    ,@(79, 55)
    ,@(80, 56)
    ,@(81, 57)
    ,@(82, 58)
    ,@(83, 59)
    ,@(84, 60)
    ,@(85, 61)
    ,@(86, 62)
    ,@(87, 63)
    ,@(88, 64)
    ,@(89, 65)
    ,@(90, 66)
    ,@(91, 89)
    ,@(92, 90)
    ,@(93, 91)
    ,@(94, 92)
    ,@(95, 93)
    ,@(96, 94)
    ,@(97, 95)
    ,@(98, 96)
    ,@(99, 97)
    ,@(100, 98)
    ,@(101, 67)
    ,@(102, 68)
    ,@(103, 69)
    ,@(104, 70)
    ,@(105, 71)
    ,@(106, 72)
    ,@(107, 73)
    ,@(108, 74)
    ,@(109, 75)
    ,@(110, 76)
    ,@(111, 77)
    ,@(112, 99)
    ,@(113, 100)
    ,@(114, 101)
    ,@(115, 102)
    ,@(116, 103)
    ,@(117, 104)
    ,@(118, 105)
    ,@(119, 106)
    ,@(120, 107)
)

foreach ($move in $moves) {
    $destRow = $move[0]
    $srcRow = $move[1]
    $ws.Cells.Item($destRow, 1).Value = $origA[$srcRow]
    $ws.Cells.Item($destRow, 2).Value = $origB[$srcRow]
}
